$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5:C9").Value = "Ruben Chavez"
$ws.Range("C5:C9").Select() | Out-Null
